# Apply the "Updated cryptos list" data refresh (GitHub Actions run).
# Only cell VALUES change (Coin/Link/Price/Volume(1h) columns);
# no rows/columns are inserted or removed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal string into a cell. Price-column values that
# look like plain numbers (e.g. "304.97") would otherwise be silently
# reinterpreted by Excel as numeric values (and re-serialized with
# floating point noise, e.g. 304.97000000000003), so such cells are
# pre-formatted as Text ("@") to force them to stay exact strings -
# matching the source workbook, where every cell is stored as text.
function Set-TextCell($addr, $value, [bool]$forceText) {
    $cell = $ws.Range($addr)
    if ($forceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $value
}

Set-TextCell "D2" "43.465.54" $False
Set-TextCell "E2" "  +0.17%  " $False
Set-TextCell "D3" "2.335.71" $False
Set-TextCell "E3" "  -1.63%  " $False
Set-TextCell "E4" "  +0.01%  " $False
Set-TextCell "D5" "304.97" $True
Set-TextCell "E5" "  -1.53%  " $False
Set-TextCell "D6" "101.47" $True
Set-TextCell "E6" "  -3.78%  " $False
Set-TextCell "D7" "0.513" $True
Set-TextCell "E7" "  -1.37%  " $False
Set-TextCell "D9" "0.515" $True
Set-TextCell "E9" "  -0.67%  " $False
Set-TextCell "D10" "35.32" $True
Set-TextCell "E10" "  -2.36%  " $False
Set-TextCell "D11" "51.65" $True
Set-TextCell "E11" "  -3.16%  " $False
Set-TextCell "D12" "0.0798" $True
Set-TextCell "E12" "  -1.85%  " $False
Set-TextCell "E13" "  +0.69%  " $False
Set-TextCell "D14" "6.84" $True
Set-TextCell "E14" "  -2.38%  " $False
Set-TextCell "D15" "2.700.58" $False
Set-TextCell "E15" "  -1.54%  " $False
Set-TextCell "D16" "15.73" $True
Set-TextCell "E16" "  +0.75%  " $False
Set-TextCell "D17" "2.342.08" $False
Set-TextCell "E17" "  -1.53%  " $False
Set-TextCell "D18" "0.807" $True
Set-TextCell "E18" "  -1.02%  " $False
Set-TextCell "D19" "43.383.96" $False
Set-TextCell "E19" "  +0.08%  " $False
Set-TextCell "D20" "11.81" $True
Set-TextCell "E20" "  -1.72%  " $False
Set-TextCell "D21" "0.0₃0908" $False
Set-TextCell "E21" "  -1.23%  " $False
Set-TextCell "D22" "6.12" $True
Set-TextCell "E22" "  -2.77%  " $False
Set-TextCell "E23" "  -0.66%  " $False
Set-TextCell "D24" "238.45" $True
Set-TextCell "E24" "  -1.41%  " $False
Set-TextCell "E25" "  -3.44%  " $False
Set-TextCell "D26" "2.54" $True
Set-TextCell "E26" "  -2.84%  " $False
Set-TextCell "D28" "25.10" $True
Set-TextCell "E28" "  -2.92%  " $False
Set-TextCell "D29" "34.78" $True
Set-TextCell "E29" "  -5.86%  " $False
Set-TextCell "B30" "Monero" $False
Set-TextCell "C30" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" $False
Set-TextCell "D30" "165.96" $True
Set-TextCell "E30" "  +2.38%  " $False
Set-TextCell "D31" "9.27" $True
Set-TextCell "E31" "  -3.41%  " $False
Set-TextCell "B32" "Toncoin" $False
Set-TextCell "C32" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton" $False
Set-TextCell "D32" "2.07" $True
Set-TextCell "E32" "  -2.29%  " $False
Set-TextCell "D33" "0.999" $True
Set-TextCell "E33" "  -0.05%  " $False
Set-TextCell "D34" "5.08" $True
Set-TextCell "E34" "  -3.66%  " $False
Set-TextCell "E35" "  -4.84%  " $False
Set-TextCell "E36" "  -6.16%  " $False
Set-TextCell "D37" "0.0710" $True
Set-TextCell "E37" "  -4.63%  " $False
Set-TextCell "D38" "16.90" $True
Set-TextCell "E38" "  -7.45%  " $False
Set-TextCell "E39" "  -6.52%  " $False
Set-TextCell "E40" "  -6.03%  " $False
Set-TextCell "D41" "0.103" $True
Set-TextCell "E41" "  -3.04%  " $False
Set-TextCell "E42" "  -2.54%  " $False
Set-TextCell "D43" "2.42" $True
Set-TextCell "E43" "  -3.11%  " $False
Set-TextCell "D44" "1.983.08" $False
Set-TextCell "D45" "0.0286" $True
Set-TextCell "E45" "  -1.51%  " $False
Set-TextCell "D46" "18.66" $True
Set-TextCell "E46" "  -6.47%  " $False
Set-TextCell "D47" "2.96" $True
Set-TextCell "E47" "  -6.28%  " $False
Set-TextCell "D48" "9.89" $True
Set-TextCell "E48" "  -6.42%  " $False
Set-TextCell "E49" "  +4.17%  " $False
Set-TextCell "D50" "55.60" $True
Set-TextCell "E50" "  -4.32%  " $False
Set-TextCell "D51" "2.562.36" $False
Set-TextCell "E51" "  +0.21%  " $False
